# Daily attendance processing - 2026-01-28 20:45:33
#
# The "Recorded By" column (G) holds a comma-separated list of the
# accounts that touched each attendance record. Normalize the ordering of
# that list for every data row: the "System" entry (exact case) always
# moves to the front (the rest keep their relative order); when "System"
# is not present in the list, the list is simply reversed.

function Get-NormalizedRecordedBy($value) {
    $items = $value -split ", "

    $hasSystem = $false
    foreach ($item in $items) {
        if ($item.CompareTo("System") -eq 0) {
            $hasSystem = $true
        }
    }

    $ordered = @()
    if ($hasSystem) {
        $ordered += "System"
        foreach ($item in $items) {
            if ($item.CompareTo("System") -ne 0) {
                $ordered += $item
            }
        }
    } else {
        for ($i = $items.Count - 1; $i -ge 0; $i--) {
            $ordered += $items[$i]
        }
    }

    return ($ordered -join ", ")
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2

    if ($current -ne $null -and $current -ne "") {
        $normalized = Get-NormalizedRecordedBy $current
        if ($normalized.CompareTo($current) -ne 0) {
            $cell.Value = $normalized
        }
    }
}
